$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date in A1 (serial date value)
$ws.Range("A1").Value = 45436

# Update prices in column D for rows 33-35
$ws.Range("D33").Value = 457
$ws.Range("D34").Value = 480
$ws.Range("D35").Value = 562
